# Actualización automática 2025-11-27 16:30:09
#
# Updates a handful of cached sales figures (and their dependent
# totals/percentages) across the three sheets of the workbook, plus the
# "N de 54" non-zero-count labels on the "VENTAS POR GRUPO" sheet and a
# minor column-width tweak on "CUMPLIMIENTO MENSUAL".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO  (per-client sales by product group)
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M4").Value  = 9030.85
$wsGrupo.Range("D24").Value = 1391.04
$wsGrupo.Range("M24").Value = 380.64
$wsGrupo.Range("D36").Value = 3810.24
$wsGrupo.Range("K36").Value = 152.28
$wsGrupo.Range("L36").Value = 3412.41
$wsGrupo.Range("M36").Value = 1717.24
$wsGrupo.Range("M37").Value = 3232.52
$wsGrupo.Range("C48").Value = 777.6
$wsGrupo.Range("F52").Value = 153.5
$wsGrupo.Range("I53").Value = 49.5
$wsGrupo.Range("M53").Value = 82.94

# Row 56 holds "<count> de 54" labels (count of non-zero rows per column).
$wsGrupo.Range("C56").Value = "4 de 54"
$wsGrupo.Range("F56").Value = "1 de 54"
$wsGrupo.Range("K56").Value = "3 de 54"
$wsGrupo.Range("M56").Value = "17 de 54"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL  (per-client monthly sales, "noviembre" column)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value  = 11751.84
$wsMensual.Range("F24").Value = 5478.92
$wsMensual.Range("F36").Value = 9092.17
$wsMensual.Range("F37").Value = 3690.44
$wsMensual.Range("F48").Value = 2139.68
$wsMensual.Range("F53").Value = 2479.21
$wsMensual.Range("F54").Value = 2479.21
$wsMensual.Range("F55").Value = 766.04
$wsMensual.Range("F56").Value = 766.04
$wsMensual.Range("F60").Value = 87549.85000000001

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL  (by product group: PRESUPUESTO/VENTA/...)
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D2").Value = 2073.6
$wsCumpl.Range("E2").Value = 753.0599999999999
$wsCumpl.Range("F2").Value = 0.7335866358175374

$wsCumpl.Range("D3").Value = 10794.53
$wsCumpl.Range("E3").Value = -4171.27
$wsCumpl.Range("F3").Value = 1.629791069654521

$wsCumpl.Range("D7").Value = 728.1
$wsCumpl.Range("E7").Value = 591.9
$wsCumpl.Range("F7").Value = 0.5515909090909091

$wsCumpl.Range("D9").Value = 161.69
$wsCumpl.Range("E9").Value = -161.69

$wsCumpl.Range("D10").Value = 976.17
$wsCumpl.Range("E10").Value = 3335.83
$wsCumpl.Range("F10").Value = 0.2263845083487941

$wsCumpl.Range("D11").Value = 20078.66
$wsCumpl.Range("E11").Value = -5842.67
$wsCumpl.Range("F11").Value = 1.410415432997635

$wsCumpl.Range("D12").Value = 45055.59
$wsCumpl.Range("E12").Value = 19888.41
$wsCumpl.Range("F12").Value = 0.6937606245380635

$wsCumpl.Range("D14").Value = 84304.59999999999
$wsCumpl.Range("E14").Value = 14651.65685923838
$wsCumpl.Range("F14").Value = 0.8519380449072581

# Column E narrows from 24 to 23 characters wide. The ColumnWidth COM
# property is offset from the raw OOXML "width" by Excel's standard
# ~0.83-character padding, so 22.17 round-trips to a stored width of 23.
$wsCumpl.Columns.Item(5).ColumnWidth = 22.17
